$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 becomes: Baggy Mavi / 350 TL / Jeans / BAG1.jpg / 100% Pamuk / Var
$ws.Range("B2").Value = "350 TL "
$ws.Range("C2").Value = "Jeans"
$ws.Range("D2").Value = "BAG1.jpg"
$ws.Range("E2").Value = "100% Pamuk"
$ws.Range("F2").Value = "Var"

# Row 3 becomes: Baggy Gri / 350 TL / Jeans / BAG2.jpg / 100% Pamuk / Var
$ws.Range("B3").Value = "350 TL "
$ws.Range("C3").Value = "Jeans"
$ws.Range("D3").Value = "BAG2.jpg"
$ws.Range("E3").Value = "100% Pamuk"
$ws.Range("F3").Value = "Var"

# Product names last, so they land at the tail of the shared-string table
$ws.Range("A2").Value = "Baggy Mavi"
$ws.Range("A3").Value = "Baggy Gri"

# The old rows 4 and 5 (elma/armut/cilek2-era rows replaced by the two Baggy
# variants above) are removed - clear them but keep row 8's placeholder
# formatting cell where it is (no entire-row delete/shift).
$ws.Range("A4:F5").ClearContents()

# Move the active selection to A4, matching the saved view state
$ws.Range("A4").Select()
